$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 132
$ws.Range("H132").Value = 4548938.5
$ws.Range("I132").Value = 6669917
$ws.Range("J132").Value = 3984.7856
$ws.Range("K132").Value = 20009751
$ws.Range("L132").Value = 11954.3568
$ws.Range("M132").Value = -20007221
$ws.Range("N132").Value = -17014.3568

# Row 137
$ws.Range("H137").Value = 2224963.2
$ws.Range("I137").Value = 3574568
$ws.Range("J137").Value = 2085.2942
$ws.Range("K137").Value = 10723704
$ws.Range("L137").Value = 6255.882599999999
$ws.Range("M137").Value = -10721154
$ws.Range("N137").Value = -11355.8826

# Row 138
$ws.Range("H138").Value = 2399.3171
$ws.Range("I138").Value = 843.46155
$ws.Range("J138").Value = 5096.1333
$ws.Range("K138").Value = 2530.38465
$ws.Range("L138").Value = 15288.3999
$ws.Range("M138").Value = 2609.61535
$ws.Range("N138").Value = -25568.3999

# Row 141
$ws.Range("H141").Value = 209599.12
$ws.Range("I141").Value = 1067.75
$ws.Range("J141").Value = 1738829.1
$ws.Range("K141").Value = 3203.25
$ws.Range("L141").Value = 5216487.300000001
$ws.Range("M141").Value = 1976.75
$ws.Range("N141").Value = -5226847.300000001

# ---------------------------------------------------------------------------
# Sheet: ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 44
$ws.Range("H44").Value = 18049
$ws.Range("J44").Value = 18049
$ws.Range("L44").Value = 18049
$ws.Range("N44").Value = -19025

# Row 51
$ws.Range("H51").Value = 80047
$ws.Range("J51").Value = 80047
$ws.Range("L51").Value = 80047
$ws.Range("N51").Value = -81559

# Row 54 - H/J/L become 0, N54 cell is removed entirely
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# Row 61
$ws.Range("H61").Value = 1718.0652
$ws.Range("I61").Value = 797.14703
$ws.Range("J61").Value = 4327.3335
$ws.Range("K61").Value = 797.14703
$ws.Range("L61").Value = 4327.3335
$ws.Range("M61").Value = -585.14703
$ws.Range("N61").Value = -4751.3335

# Row 74
$ws.Range("H74").Value = 809.04
$ws.Range("I74").Value = 705.5217
$ws.Range("J74").Value = 1999.5
$ws.Range("K74").Value = 705.5217
$ws.Range("L74").Value = 1999.5
$ws.Range("M74").Value = 168.4783
$ws.Range("N74").Value = -3747.5

# Row 77
$ws.Range("H77").Value = 809.04
$ws.Range("I77").Value = 705.5217
$ws.Range("J77").Value = 1999.5
$ws.Range("K77").Value = 3527.6085
$ws.Range("L77").Value = 9997.5
$ws.Range("M77").Value = 840.3914999999997
$ws.Range("N77").Value = -18733.5

# Row 136
$ws.Range("H136").Value = 1718.0652
$ws.Range("I136").Value = 797.14703
$ws.Range("J136").Value = 4327.3335
$ws.Range("K136").Value = 2391.44109
$ws.Range("L136").Value = 12982.0005
$ws.Range("M136").Value = 158.5589100000002
$ws.Range("N136").Value = -18082.0005

# ---------------------------------------------------------------------------
# Sheet: BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 22
$ws.Range("H22").Value = 425
$ws.Range("I22").Value = 360
$ws.Range("J22").Value = 490
$ws.Range("K22").Value = 360
$ws.Range("L22").Value = 490
$ws.Range("M22").Value = -187
$ws.Range("N22").Value = -836

# Row 134
$ws.Range("H134").Value = 1715.762
$ws.Range("I134").Value = 823.94446
$ws.Range("J134").Value = 7066.6665
$ws.Range("K134").Value = 2471.83338
$ws.Range("L134").Value = 21199.9995
$ws.Range("M134").Value = 63.16661999999997
$ws.Range("N134").Value = -26269.9995

# Row 137 - J/L change, N137 is newly added, M137 unchanged
$ws.Range("H137").Value = 33578.43
$ws.Range("J137").Value = 33890
$ws.Range("L137").Value = 33890
$ws.Range("N137").Value = -44090

# ---------------------------------------------------------------------------
# Sheet: CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 3452394.2
$ws.Range("I31").Value = 5884761
$ws.Range("J31").Value = 6541.1665
$ws.Range("K31").Value = 5884761
$ws.Range("L31").Value = 6541.1665
$ws.Range("M31").Value = -5884466
$ws.Range("N31").Value = -7131.1665

# Row 34
$ws.Range("H34").Value = 3452394.2
$ws.Range("I34").Value = 5884761
$ws.Range("J34").Value = 6541.1665
$ws.Range("K34").Value = 5884761
$ws.Range("L34").Value = 6541.1665
$ws.Range("M34").Value = -5884559
$ws.Range("N34").Value = -6945.1665

# Row 132
$ws.Range("H132").Value = 2860.8
$ws.Range("I132").Value = 2128.4285
$ws.Range("J132").Value = 4569.6665
$ws.Range("K132").Value = 6385.2855
$ws.Range("L132").Value = 13708.9995
$ws.Range("M132").Value = -3855.2855
$ws.Range("N132").Value = -18768.9995

# ---------------------------------------------------------------------------
# Sheet: CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 3
$ws.Range("H3").Value = 3547.6086
$ws.Range("I3").Value = 1635
$ws.Range("J3").Value = 5018.846
$ws.Range("K3").Value = 4905
$ws.Range("L3").Value = 15056.538
$ws.Range("M3").Value = -4793
$ws.Range("N3").Value = -15280.538

# Row 5
$ws.Range("H5").Value = 2045.3334
$ws.Range("I5").Value = 1302
$ws.Range("J5").Value = 2640
$ws.Range("K5").Value = 3906
$ws.Range("L5").Value = 7920
$ws.Range("M5").Value = -3794
$ws.Range("N5").Value = -8144

# Row 135
$ws.Range("H135").Value = 2045.3334
$ws.Range("I135").Value = 1302
$ws.Range("J135").Value = 2640
$ws.Range("K135").Value = 11718
$ws.Range("L135").Value = 23760
$ws.Range("M135").Value = -9183
$ws.Range("N135").Value = -28830

# ---------------------------------------------------------------------------
# Sheet: GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 95
$ws.Range("H95").Value = 19122
$ws.Range("J95").Value = 19122
$ws.Range("L95").Value = 19122
$ws.Range("N95").Value = -24614

# Row 122
$ws.Range("H122").Value = 4022.0293
$ws.Range("I122").Value = 3173.2666
$ws.Range("J122").Value = 4692.1055
$ws.Range("K122").Value = 9519.799800000001
$ws.Range("L122").Value = 14076.3165
$ws.Range("M122").Value = -7069.799800000001
$ws.Range("N122").Value = -18976.3165

# Row 132
$ws.Range("H132").Value = 3060.6155
$ws.Range("I132").Value = 2903.875
$ws.Range("J132").Value = 3311.4
$ws.Range("K132").Value = 8711.625
$ws.Range("L132").Value = 9934.200000000001
$ws.Range("M132").Value = -6181.625
$ws.Range("N132").Value = -14994.2

# ---------------------------------------------------------------------------
# Sheet: LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 24
$ws.Range("H24").Value = 25007
$ws.Range("J24").Value = 25007
$ws.Range("L24").Value = 25007
$ws.Range("N24").Value = -25693

# Row 45 - H/I/K change, J/L/N unchanged, M45 newly added
$ws.Range("H45").Value = 8199.799999999999
$ws.Range("I45").Value = 6999.5
$ws.Range("K45").Value = 6999.5
$ws.Range("M45").Value = -6592.5

# Row 132
$ws.Range("H132").Value = 3056.2
$ws.Range("I132").Value = 2086.4614
$ws.Range("J132").Value = 4857.143
$ws.Range("K132").Value = 6259.3842
$ws.Range("L132").Value = 14571.429
$ws.Range("M132").Value = -3729.3842
$ws.Range("N132").Value = -19631.429

# ---------------------------------------------------------------------------
# Sheet: WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 109
$ws.Range("H109").Value = 37515.4
$ws.Range("J109").Value = 37515.4
$ws.Range("L109").Value = 37515.4
$ws.Range("N109").Value = -40289.4

# Row 136
$ws.Range("H136").Value = 973.7447
$ws.Range("I136").Value = 473.9
$ws.Range("J136").Value = 1855.8235
$ws.Range("K136").Value = 1421.7
$ws.Range("L136").Value = 5567.470499999999
$ws.Range("M136").Value = 1128.3
$ws.Range("N136").Value = -10667.4705

Write-Host "All updates applied"
